# month_3_wastes.xlsx -> add the 5th week ("01.05-07.05 2017", row 5) figures
# to the "Lunches" sheet, and leave the workbook with "Lunches" as the
# active/selected sheet (cell H6 selected), matching the author's manual
# edit session. Downstream formulas (Lunches!I5, Lunches!I6, Total!B2,
# Total!B7, ...) recalculate automatically.

$wb = $excel.ActiveWorkbook

$wsLunches = $wb.Worksheets.Item("Lunches")

# Fill in the previously-empty Tuesday..Friday waste figures for the
# "01.05-07.05 2017" week (row 5). Column I (Total for week) is a formula
# and recalculates on its own.
$wsLunches.Range("E5").Value = 0
$wsLunches.Range("F5").Value = 0
$wsLunches.Range("G5").Value = 5.3
$wsLunches.Range("H5").Value = 0

# Make "Lunches" the active sheet and leave the selection on H6, as left by
# whoever made the edit.
$wsLunches.Activate() | Out-Null
$wsLunches.Range("H6").Select() | Out-Null
